$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.672.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.818.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.58"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.562"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.21%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.078.18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.42"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.811.57"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.18%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.686.21"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.28"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.66"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0803"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.74%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "171.45"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.66%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.64%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.37%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.30%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.87%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.420.23"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "86.42"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.02%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.47%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.45%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.86"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.20%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.11"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.980.62"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.14"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.02%  "
